# Apply the cibmtr-reporting-ig update to the "med-carmustine" ValueSet workbook.
#
# Sheet "Metadata" (sheet1) property/value table changes:
#   - Version      0.1.6 -> 0.1.7
#   - Status       active -> draft
#   - Date         2023-05-05T10:50:04-05:00 -> 2024-08-27T12:23:18-05:00
#   - Contact #1   "No display for ContactDetail" ->
#                  "The Medical College of Wisconsin, Inc. and the National
#                   Marrow Donor Program (http://www.cibmtr.org)"
#   - Contact #2   "No display for ContactDetail" ->
#                  "Bob Milius (bmilius@nmdp.org)"
#   - a new "Jurisdiction" / "" row is inserted right after the Contact rows,
#     pushing Description/Purpose/Copyright/Immutable down by one row.
#
# Sheet "Include from RxNorm" (sheet2) has no content changes.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- simple text/value updates -------------------------------------------------
$ws1.Range("B3").Value = "0.1.7"
$ws1.Range("B6").Value = "draft"
$ws1.Range("B8").Value = "2024-08-27T12:23:18-05:00"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- insert a new "Jurisdiction" row after the two Contact rows ---------------
# Shift existing rows 12-15 down to 13-16, working bottom-up so a source row is
# never overwritten before it has been copied. Copy(Destination) carries over
# both value and style (border/wrap formatting) without forking new styles.
$ws1.Range("A15:B15").Copy($ws1.Range("A16:B16"))

$ws1.Range("A14:B14").Copy($ws1.Range("A15:B15"))
$ws1.Range("B15").ClearContents()

$ws1.Range("A13:B13").Copy($ws1.Range("A14:B14"))
$ws1.Range("B14").ClearContents()

$ws1.Range("A12:B12").Copy($ws1.Range("A13:B13"))

# Row 12 becomes the new Jurisdiction row (value is empty, same style as the
# other data rows, inherited from the copy above).
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").ClearContents()
